$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 8.725601455334775
$ws.Cells.Item(2, 3).Value = 5.781599240996697
$ws.Cells.Item(2, 4).Value = 5.980611105102169
$ws.Cells.Item(2, 5).Value = 12.76293396798647
$ws.Cells.Item(2, 7).Value = 28.10905514645271
$ws.Cells.Item(2, 8).Value = 14.09959807047985
$ws.Cells.Item(2, 11).Value = 7.936820453847593
$ws.Cells.Item(2, 13).Value = 13.37237366917323
$ws.Cells.Item(2, 14).Value = 18.43650546021523
$ws.Cells.Item(2, 15).Value = 21.38109274803454
$ws.Cells.Item(3, 2).Value = 8.436408262567735
$ws.Cells.Item(3, 3).Value = 5.675404925169916
$ws.Cells.Item(3, 4).Value = 5.860764096022106
$ws.Cells.Item(3, 5).Value = 12.54969614024091
$ws.Cells.Item(3, 7).Value = 28.10051850850035
$ws.Cells.Item(3, 8).Value = 14.13792960440271
$ws.Cells.Item(3, 11).Value = 7.740783164299862
$ws.Cells.Item(3, 13).Value = 13.19933193705617
$ws.Cells.Item(3, 14).Value = 18.49610281837917
$ws.Cells.Item(3, 15).Value = 21.43055067625185
$ws.Cells.Item(4, 2).Value = 8.255269698877035
$ws.Cells.Item(4, 3).Value = 5.60850043691454
$ws.Cells.Item(4, 4).Value = 5.787659500711089
$ws.Cells.Item(4, 5).Value = 12.42134070929188
$ws.Cells.Item(4, 7).Value = 28.10421294811156
$ws.Cells.Item(4, 8).Value = 14.16369894441982
$ws.Cells.Item(4, 11).Value = 7.619158236679636
$ws.Cells.Item(4, 13).Value = 13.09521608548949
$ws.Cells.Item(4, 14).Value = 18.53438309409822
$ws.Cells.Item(4, 15).Value = 21.46549919259889
$ws.Cells.Item(5, 2).Value = 8.180676287135439
$ws.Cells.Item(5, 3).Value = 5.580829208160272
$ws.Cells.Item(5, 4).Value = 5.758036618526384
$ws.Cells.Item(5, 5).Value = 12.36975729847239
$ws.Cells.Item(5, 7).Value = 28.10796369685828
$ws.Cells.Item(5, 8).Value = 14.1747615566878
$ws.Cells.Item(5, 11).Value = 7.569353132714046
$ws.Cells.Item(5, 13).Value = 13.05337536051415
$ws.Cells.Item(5, 14).Value = 18.55040815924673
$ws.Cells.Item(5, 15).Value = 21.48089063792745
$ws.Cells.Item(6, 2).Value = 8.168246831190698
$ws.Cells.Item(6, 3).Value = 5.576210451798395
$ws.Cells.Item(6, 4).Value = 5.753129271597144
$ws.Cells.Item(6, 5).Value = 12.36123772553355
$ws.Cells.Item(6, 7).Value = 28.10872202279387
$ws.Cells.Item(6, 8).Value = 14.17663239393658
$ws.Cells.Item(6, 11).Value = 7.561070652270621
$ws.Cells.Item(6, 13).Value = 13.04646468462045
$ws.Cells.Item(6, 14).Value = 18.55309484807864
$ws.Cells.Item(6, 15).Value = 21.48351573922007
$ws.Cells.Item(7, 2).Value = 8.254266689710308
$ws.Cells.Item(7, 3).Value = 5.608128873583142
$ws.Cells.Item(7, 4).Value = 5.787259254574809
$ws.Cells.Item(7, 5).Value = 12.42064201346161
$ws.Cells.Item(7, 7).Value = 28.10425444496226
$ws.Cells.Item(7, 8).Value = 14.16384586597534
$ws.Cells.Item(7, 11).Value = 7.618487426846535
$ws.Cells.Item(7, 13).Value = 13.09464936207989
$ws.Cells.Item(7, 14).Value = 18.5345974890533
$ws.Cells.Item(7, 15).Value = 21.46570211540201
$ws.Cells.Item(8, 2).Value = 8.626698375533124
$ws.Cells.Item(8, 3).Value = 5.745348481481644
$ws.Cells.Item(8, 4).Value = 5.939215613333714
$ws.Cells.Item(8, 5).Value = 12.68892097149961
$ws.Cells.Item(8, 7).Value = 28.10425654818285
$ws.Cells.Item(8, 8).Value = 14.11235104065991
$ws.Cells.Item(8, 11).Value = 7.869531180324733
$ws.Cells.Item(8, 13).Value = 13.31229725808578
$ws.Cells.Item(8, 14).Value = 18.45670520598119
$ws.Cells.Item(8, 15).Value = 21.39719368448924
$ws.Cells.Item(9, 2).Value = 9.323899034017421
$ws.Cells.Item(9, 3).Value = 6.000086513129707
$ws.Cells.Item(9, 4).Value = 6.239016172089172
$ws.Cells.Item(9, 5).Value = 13.23208559159682
$ws.Cells.Item(9, 7).Value = 28.17515806936498
$ws.Cells.Item(9, 8).Value = 14.02910418788099
$ws.Cells.Item(9, 11).Value = 8.413751216739517
$ws.Cells.Item(9, 13).Value = 13.75373591995683
$ws.Cells.Item(9, 14).Value = 18.317287354759
$ws.Cells.Item(9, 15).Value = 21.29929489277364
$ws.Cells.Item(10, 2).Value = 9.810219128045244
$ws.Cells.Item(10, 3).Value = 6.177422010430123
$ws.Cells.Item(10, 4).Value = 6.457781602760402
$ws.Cells.Item(10, 5).Value = 13.63707703587992
$ws.Cells.Item(10, 7).Value = 28.27033850416893
$ws.Cells.Item(10, 8).Value = 13.97876929792985
$ws.Cells.Item(10, 11).Value = 8.87393216611496
$ws.Cells.Item(10, 13).Value = 14.08388626544473
$ws.Cells.Item(10, 14).Value = 18.22289849158818
$ws.Cells.Item(10, 15).Value = 21.24970821123426
$ws.Cells.Item(11, 2).Value = 10.02481811808395
$ws.Cells.Item(11, 3).Value = 6.255757442531566
$ws.Cells.Item(11, 4).Value = 6.556460410622266
$ws.Cells.Item(11, 5).Value = 13.82167681742981
$ws.Cells.Item(11, 7).Value = 28.32292448578382
$ws.Cells.Item(11, 8).Value = 13.95822447464865
$ws.Cells.Item(11, 11).Value = 9.073976219844202
$ws.Cells.Item(11, 13).Value = 14.23469047135316
$ws.Cells.Item(11, 14).Value = 18.18168670461345
$ws.Cells.Item(11, 15).Value = 21.23202198067202
$ws.Cells.Item(12, 2).Value = 10.10505548135608
$ws.Cells.Item(12, 3).Value = 6.28506927985956
$ws.Cells.Item(12, 4).Value = 6.593666567675059
$ws.Cells.Item(12, 5).Value = 13.8915580264062
$ws.Cells.Item(12, 7).Value = 28.34416371588687
$ws.Cells.Item(12, 8).Value = 13.9507832329556
$ws.Cells.Item(12, 11).Value = 9.148360144919229
$ws.Cells.Item(12, 13).Value = 14.29183106401439
$ws.Cells.Item(12, 14).Value = 18.16632777059509
$ws.Cells.Item(12, 15).Value = 21.22602640110147
$ws.Cells.Item(13, 2).Value = 10.08782166426302
$ws.Cells.Item(13, 3).Value = 6.27877235825897
$ws.Cells.Item(13, 4).Value = 6.585661357099984
$ws.Cells.Item(13, 5).Value = 13.87651004148009
$ws.Cells.Item(13, 7).Value = 28.33953066431714
$ws.Cells.Item(13, 8).Value = 13.95237077184443
$ws.Cells.Item(13, 11).Value = 9.132401554331308
$ws.Cells.Item(13, 13).Value = 14.27952414292283
$ws.Cells.Item(13, 14).Value = 18.16962462126352
$ws.Cells.Item(13, 15).Value = 21.22728642428166
$ws.Cells.Item(14, 2).Value = 10.03144025833098
$ws.Cells.Item(14, 3).Value = 6.258176090843098
$ws.Cells.Item(14, 4).Value = 6.559524842780351
$ws.Cells.Item(14, 5).Value = 13.82742686476686
$ws.Cells.Item(14, 7).Value = 28.32464533568132
$ws.Cells.Item(14, 8).Value = 13.95760549073677
$ws.Cells.Item(14, 11).Value = 9.080123421151225
$ws.Cells.Item(14, 13).Value = 14.23939108219135
$ws.Cells.Item(14, 14).Value = 18.18041817064782
$ws.Cells.Item(14, 15).Value = 21.23151464850927
$ws.Cells.Item(15, 2).Value = 9.996769296098481
$ws.Cells.Item(15, 3).Value = 6.245513965136034
$ws.Cells.Item(15, 4).Value = 6.543493288325031
$ws.Cells.Item(15, 5).Value = 13.79735684782017
$ws.Cells.Item(15, 7).Value = 28.3157000163969
$ws.Cells.Item(15, 8).Value = 13.96085601653399
$ws.Cells.Item(15, 11).Value = 9.047922481200587
$ws.Cells.Item(15, 13).Value = 14.21481129673739
$ws.Cells.Item(15, 14).Value = 18.1870616739187
$ws.Cells.Item(15, 15).Value = 21.23419599052797
$ws.Cells.Item(16, 2).Value = 9.796054708735879
$ws.Cells.Item(16, 3).Value = 6.172254197176015
$ws.Cells.Item(16, 4).Value = 6.451312296857297
$ws.Cells.Item(16, 5).Value = 13.62501387789531
$ws.Cells.Item(16, 7).Value = 28.26708801544562
$ws.Cells.Item(16, 8).Value = 13.98015932142826
$ws.Cells.Item(16, 11).Value = 8.860668953281122
$ws.Cells.Item(16, 13).Value = 14.07403859572887
$ws.Cells.Item(16, 14).Value = 18.22562641564167
$ws.Cells.Item(16, 15).Value = 21.25096218814293
$ws.Cells.Item(17, 2).Value = 9.67116980986869
$ws.Cells.Item(17, 3).Value = 6.126701645094249
$ws.Cells.Item(17, 4).Value = 6.394517706691953
$ws.Cells.Item(17, 5).Value = 13.51932652586571
$ws.Cells.Item(17, 7).Value = 28.23963885878987
$ws.Cells.Item(17, 8).Value = 13.99260410461115
$ws.Cells.Item(17, 11).Value = 8.743390212657909
$ws.Cells.Item(17, 13).Value = 13.98779809468738
$ws.Cells.Item(17, 14).Value = 18.24972593224011
$ws.Cells.Item(17, 15).Value = 21.2624963468777
$ws.Cells.Item(18, 2).Value = 9.598719040367341
$ws.Cells.Item(18, 3).Value = 6.100282110326568
$ws.Cells.Item(18, 4).Value = 6.361774316483286
$ws.Cells.Item(18, 5).Value = 13.458576878933
$ws.Cells.Item(18, 7).Value = 28.22472575515657
$ws.Cells.Item(18, 8).Value = 13.99998345160321
$ws.Cells.Item(18, 11).Value = 8.675061317089645
$ws.Cells.Item(18, 13).Value = 13.93825612504971
$ws.Cells.Item(18, 14).Value = 18.26374988095339
$ws.Cells.Item(18, 15).Value = 21.26958890655853
$ws.Cells.Item(19, 2).Value = 9.574084290638107
$ws.Cells.Item(19, 3).Value = 6.091299824601268
$ws.Cells.Item(19, 4).Value = 6.350676083193194
$ws.Cells.Item(19, 5).Value = 13.43801712133501
$ws.Cells.Item(19, 7).Value = 28.21982696056774
$ws.Cells.Item(19, 8).Value = 14.0025199937914
$ws.Cells.Item(19, 11).Value = 8.651777392370883
$ws.Cells.Item(19, 13).Value = 13.92149426806867
$ws.Cells.Item(19, 14).Value = 18.26852610128028
$ws.Cells.Item(19, 15).Value = 21.2720690127823
$ws.Cells.Item(20, 2).Value = 9.684528773987628
$ws.Cells.Item(20, 3).Value = 6.131573569768703
$ws.Cells.Item(20, 4).Value = 6.400571822782187
$ws.Cells.Item(20, 5).Value = 13.53057365209039
$ws.Cells.Item(20, 7).Value = 28.24247037939865
$ws.Cells.Item(20, 8).Value = 13.99125641596358
$ws.Cells.Item(20, 11).Value = 8.75596536039936
$ws.Cells.Item(20, 13).Value = 13.99697259835567
$ws.Cells.Item(20, 14).Value = 18.24714368317634
$ws.Cells.Item(20, 15).Value = 21.26122106079432
$ws.Cells.Item(21, 2).Value = 10.04802924738365
$ws.Cells.Item(21, 3).Value = 6.264235393846778
$ws.Cells.Item(21, 4).Value = 6.567206463341131
$ws.Cells.Item(21, 5).Value = 13.84184498850991
$ws.Cells.Item(21, 7).Value = 28.32898161474754
$ws.Cells.Item(21, 8).Value = 13.95605873404407
$ws.Cells.Item(21, 11).Value = 9.095516145326995
$ws.Cells.Item(21, 13).Value = 14.25117863165137
$ws.Cells.Item(21, 14).Value = 18.17724114709519
$ws.Cells.Item(21, 15).Value = 21.23025366003331
$ws.Cells.Item(22, 2).Value = 10.27958578907657
$ws.Cells.Item(22, 3).Value = 6.348878513634759
$ws.Cells.Item(22, 4).Value = 6.675154027711804
$ws.Cells.Item(22, 5).Value = 14.04511334494828
$ws.Cells.Item(22, 7).Value = 28.39324514309727
$ws.Cells.Item(22, 8).Value = 13.935028975371
$ws.Cells.Item(22, 11).Value = 9.309444446179576
$ws.Cells.Item(22, 13).Value = 14.41749262005115
$ws.Cells.Item(22, 14).Value = 18.13299541744779
$ws.Cells.Item(22, 15).Value = 21.21410581125961
$ws.Cells.Item(23, 2).Value = 10.15657182916806
$ws.Cells.Item(23, 3).Value = 6.303896390487622
$ws.Cells.Item(23, 4).Value = 6.617640880736186
$ws.Cells.Item(23, 5).Value = 13.93666456023246
$ws.Cells.Item(23, 7).Value = 28.35824346985288
$ws.Cells.Item(23, 8).Value = 13.94607223978249
$ws.Cells.Item(23, 11).Value = 9.196007109557582
$ws.Cells.Item(23, 13).Value = 14.32872930521067
$ws.Cells.Item(23, 14).Value = 18.15647884736507
$ws.Cells.Item(23, 15).Value = 21.22234950178068
$ws.Cells.Item(24, 2).Value = 9.678491213915844
$ws.Cells.Item(24, 3).Value = 6.129371690185085
$ws.Cells.Item(24, 4).Value = 6.397835038548674
$ws.Cells.Item(24, 5).Value = 13.5254887849773
$ws.Cells.Item(24, 7).Value = 28.24118754512934
$ws.Cells.Item(24, 8).Value = 13.99186500616114
$ws.Cells.Item(24, 11).Value = 8.750282948767056
$ws.Cells.Item(24, 13).Value = 13.99282468008543
$ws.Cells.Item(24, 14).Value = 18.24831059202186
$ws.Cells.Item(24, 15).Value = 21.26179618043134
$ws.Cells.Item(25, 2).Value = 9.139463219636912
$ws.Cells.Item(25, 3).Value = 5.932824024010191
$ws.Cells.Item(25, 4).Value = 6.1579928936127
$ws.Cells.Item(25, 5).Value = 13.0837830997991
$ws.Cells.Item(25, 7).Value = 28.14839283379655
$ws.Cells.Item(25, 8).Value = 14.04972442792112
$ws.Cells.Item(25, 11).Value = 8.236283570269817
$ws.Cells.Item(25, 13).Value = 13.63306800534306
$ws.Cells.Item(25, 14).Value = 18.35358537977885
$ws.Cells.Item(25, 15).Value = 21.32186387300137

Write-Output "Done updating loading_percent values"